# Updated cryptos list on Fri Apr 26 19:45:24 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force text storage (avoid numeric auto-conversion) without leaving a
    # lingering custom style on the cell: prefix with an apostrophe (Excel
    # "treat as text" marker), then reset the cell style back to Normal.
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "63.920.17"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "3.144.35"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue $ws.Range("D5") "603.64"
$ws.Range("E5").Value = "  -2.06%  "

Set-TextValue $ws.Range("D6") "143.25"
$ws.Range("E6").Value = "  -2.40%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.138.58"
$ws.Range("E8").Value = "  -0.46%  "

Set-TextValue $ws.Range("D9") "0.527"
$ws.Range("E9").Value = "  -0.48%  "

$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("E12").Value = "  -1.42%  "

Set-TextValue $ws.Range("D13") "0.0000253"
$ws.Range("E13").Value = "  -2.05%  "

Set-TextValue $ws.Range("D14") "34.89"
$ws.Range("E14").Value = "  -2.59%  "

$ws.Range("D15").Value = "3.663.35"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("E16").Value = "  +2.92%  "

$ws.Range("D17").Value = "63.938.48"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").Value = "3.151.90"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("E19").Value = "  -0.95%  "

Set-TextValue $ws.Range("D20") "487.51"
$ws.Range("E20").Value = "  +1.95%  "

Set-TextValue $ws.Range("D21") "14.69"
$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("E22").Value = "  -1.34%  "

$ws.Range("E23").Value = "  -3.09%  "

Set-TextValue $ws.Range("D24") "88.24"
$ws.Range("E24").Value = "  +4.59%  "

Set-TextValue $ws.Range("D25") "13.30"
$ws.Range("E25").Value = "  -3.38%  "

$ws.Range("E27").Value = "  -2.36%  "

Set-TextValue $ws.Range("D28") "8.20"
$ws.Range("E28").Value = "  -4.18%  "

$ws.Range("E29").Value = "  +1.08%  "

$ws.Range("E30").Value = "  -1.11%  "

Set-TextValue $ws.Range("D31") "27.55"
$ws.Range("E31").Value = "  +3.77%  "

$ws.Range("E32").Value = "  -6.13%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  -1.86%  "

$ws.Range("E35").Value = "  -2.68%  "

Set-TextValue $ws.Range("D36") "6.04"
$ws.Range("E36").Value = "  +0.52%  "

Set-TextValue $ws.Range("D37") "52.74"
$ws.Range("E37").Value = "  -0.61%  "

$ws.Range("D38").Value = "0.0₃0741"
$ws.Range("E38").Value = "  -4.90%  "

Set-TextValue $ws.Range("D39") "2.97"
$ws.Range("E39").Value = "  -6.31%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D40") "433.94"
$ws.Range("E40").Value = "  -5.70%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D41") "0.0397"
$ws.Range("E41").Value = "  -0.63%  "

$ws.Range("E42").Value = "  -0.19%  "

Set-TextValue $ws.Range("D43") "8.32"
$ws.Range("E43").Value = "  -0.89%  "

$ws.Range("D44").Value = "2.931.03"

$ws.Range("E45").Value = "  -2.92%  "

$ws.Range("E46").Value = "  -5.31%  "

$ws.Range("E47").Value = "  -0.70%  "

$ws.Range("E48").Value = "  -0.10%  "

Set-TextValue $ws.Range("D49") "25.85"
$ws.Range("E49").Value = "  -2.66%  "

$ws.Range("E50").Value = "  +0.14%  "

Set-TextValue $ws.Range("D51") "120.41"
$ws.Range("E51").Value = "  +0.12%  "

